# updates database name to refer to ecoinvent 3.8 cut-off
$wb = $excel.ActiveWorkbook

$newName = "ecoinvent 3.8 cut-off"

# "Copy Activities" sheet: source_database column (A2:A6) held the old name
$wsCopy = $wb.Worksheets.Item("Copy Activities")
$wsCopy.Activate() | Out-Null
$wsCopy.Range("A2:A6").Value = $newName
$wsCopy.Range("D2").Select() | Out-Null

# "Add Exchanges" sheet: exchange_database cell (B7) held the old name
$wsAdd = $wb.Worksheets.Item("Add Exchanges")
$wsAdd.Activate() | Out-Null
$wsAdd.Range("B7").Value = $newName
$wsAdd.Range("B7").Select() | Out-Null

# "Delete Exchanges" sheet: exchange_database column (D2:D5) held the old name
$wsDelete = $wb.Worksheets.Item("Delete Exchanges")
$wsDelete.Activate() | Out-Null
$wsDelete.Range("D2:D5").Value = $newName
$wsDelete.Range("C12").Select() | Out-Null
